# The author was "getting to know the code" - exploring the workbook and
# adding a new entity ("StorageBuckets") to the EntityNames sheet, while
# also poking at cell formatting (an underline applied to an empty cell)
# before finally leaving the selection on a couple of cells.

$wb = $excel.ActiveWorkbook

# --- EntityNames sheet -----------------------------------------------
$ws1 = $wb.Worksheets.Item("EntityNames")

# New entity added at the bottom of the list.
$ws1.Range("A8").Value = "StorageBuckets"

# Some exploratory formatting: underline applied to an (otherwise empty) cell.
$ws1.Range("F3").Font.Underline = 1

# Page setup touched for this sheet.
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# --- EntityMapping sheet ----------------------------------------------
$ws2 = $wb.Worksheets.Item("EntityMapping")
[void]$ws2.Range("A9").Select()

# Leave the cursor/selection where the author left off on EntityNames,
# which is also the active sheet.
[void]$ws1.Select()
[void]$ws1.Range("F3").Select()
